$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.0522727272727273
$ws.Cells.Item(2, 3).Value = 0.00833333333333333
$ws.Cells.Item(2, 4).Value = 0.00984848484848485
$ws.Cells.Item(2, 5).Value = 0.0128787878787879
$ws.Cells.Item(2, 6).Value = 0.940909090909091
$ws.Cells.Item(2, 7).Value = 0.0121212121212121
$ws.Cells.Item(2, 8).Value = 0.977272727272727
$ws.Cells.Item(2, 9).Value = 0.963636363636364
$ws.Cells.Item(2, 10).Value = 0.0143939393939394
$ws.Cells.Item(2, 11).Value = 0.00606060606060606
$ws.Cells.Item(2, 12).Value = 0.0053030303030303
$ws.Cells.Item(2, 13).Value = 0.992424242424242
$ws.Cells.Item(2, 14).Value = 0
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = 0.00378787878787879
$ws.Cells.Item(2, 17).Value = 0
$ws.Cells.Item(2, 18).Value = 0
$ws.Cells.Item(2, 19).Value = 0
$ws.Cells.Item(2, 20).Value = 0.0378787878787879
$ws.Cells.Item(2, 21).Value = 0.0295454545454545
$ws.Cells.Item(2, 22).Value = 0.00303030303030303
$ws.Cells.Item(2, 23).Value = 0.0143939393939394
$ws.Cells.Item(2, 24).Value = 0.0166666666666667

$ws.Cells.Item(3, 2).Value = 0.0113636363636364
$ws.Cells.Item(3, 3).Value = 0.0143939393939394
$ws.Cells.Item(3, 4).Value = 0.964393939393939
$ws.Cells.Item(3, 5).Value = 0.0136363636363636
$ws.Cells.Item(3, 6).Value = 0.00757575757575758
$ws.Cells.Item(3, 7).Value = 0.971212121212121
$ws.Cells.Item(3, 8).Value = 0.0128787878787879
$ws.Cells.Item(3, 9).Value = 0.00227272727272727
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0.00227272727272727
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = 0
$ws.Cells.Item(3, 14).Value = 0.998484848484849
$ws.Cells.Item(3, 15).Value = 0.00227272727272727
$ws.Cells.Item(3, 16).Value = 0.000757575757575758
$ws.Cells.Item(3, 17).Value = 0.00227272727272727
$ws.Cells.Item(3, 18).Value = 0.993939393939394
$ws.Cells.Item(3, 19).Value = 0.99469696969697
$ws.Cells.Item(3, 20).Value = 0
$ws.Cells.Item(3, 21).Value = 0.953787878787879
$ws.Cells.Item(3, 22).Value = 0.0234848484848485
$ws.Cells.Item(3, 23).Value = 0.00227272727272727
$ws.Cells.Item(3, 24).Value = 0

$ws.Cells.Item(4, 2).Value = 0.00984848484848485
$ws.Cells.Item(4, 3).Value = 0.187878787878788
$ws.Cells.Item(4, 4).Value = 0.00681818181818182
$ws.Cells.Item(4, 5).Value = 0.00909090909090909
$ws.Cells.Item(4, 6).Value = 0.0462121212121212
$ws.Cells.Item(4, 7).Value = 0.00227272727272727
$ws.Cells.Item(4, 8).Value = 0.0053030303030303
$ws.Cells.Item(4, 9).Value = 0.0295454545454545
$ws.Cells.Item(4, 10).Value = 0.984848484848485
$ws.Cells.Item(4, 11).Value = 0.991666666666667
$ws.Cells.Item(4, 12).Value = 0.992424242424242
$ws.Cells.Item(4, 13).Value = 0.00757575757575758
$ws.Cells.Item(4, 14).Value = 0
$ws.Cells.Item(4, 15).Value = 0.000757575757575758
$ws.Cells.Item(4, 16).Value = 0
$ws.Cells.Item(4, 17).Value = 0
$ws.Cells.Item(4, 18).Value = 0
$ws.Cells.Item(4, 19).Value = 0.000757575757575758
$ws.Cells.Item(4, 20).Value = 0.961363636363636
$ws.Cells.Item(4, 21).Value = 0.00454545454545455
$ws.Cells.Item(4, 22).Value = 0.00454545454545455
$ws.Cells.Item(4, 23).Value = 0.975757575757576
$ws.Cells.Item(4, 24).Value = 0.981818181818182

$ws.Cells.Item(5, 2).Value = 0.926515151515152
$ws.Cells.Item(5, 3).Value = 0.789393939393939
$ws.Cells.Item(5, 4).Value = 0.0189393939393939
$ws.Cells.Item(5, 5).Value = 0.964393939393939
$ws.Cells.Item(5, 6).Value = 0.0053030303030303
$ws.Cells.Item(5, 7).Value = 0.0143939393939394
$ws.Cells.Item(5, 8).Value = 0.00454545454545455
$ws.Cells.Item(5, 9).Value = 0.00454545454545455
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 0.00227272727272727
$ws.Cells.Item(5, 13).Value = 0
$ws.Cells.Item(5, 14).Value = 0.000757575757575758
$ws.Cells.Item(5, 15).Value = 0.996969696969697
$ws.Cells.Item(5, 16).Value = 0.995454545454545
$ws.Cells.Item(5, 17).Value = 0.997727272727273
$ws.Cells.Item(5, 18).Value = 0.00606060606060606
$ws.Cells.Item(5, 19).Value = 0.00454545454545455
$ws.Cells.Item(5, 20).Value = 0.000757575757575758
$ws.Cells.Item(5, 21).Value = 0.0121212121212121
$ws.Cells.Item(5, 22).Value = 0.968939393939394
$ws.Cells.Item(5, 23).Value = 0.00681818181818182
$ws.Cells.Item(5, 24).Value = 0.000757575757575758
